$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2018 LEAVE CREDITS")

# ---------------------------------------------------------------------------
# Insert a new row inside the Table15 date sequence (right after 12/31/2023,
# absolute sheet row 36) to hold a "2024" year-separator label. This pushes
# every following row (including the table's closing/last styled row) down
# by one, growing the table from A8:K83 to A8:K84.
# ---------------------------------------------------------------------------
$ws.Rows.Item(37).Insert(-4121, 0)

$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A8:K84"))

# The freshly inserted row 37 starts out with default/blank formatting.
# Clone the normal data-row formatting (which row 38 still has - it used to
# be row 37 before the shift) onto it.
$ws.Range("A38:K38").Copy()
$ws.Range("A37:K37").PasteSpecial(-4122)

# A format-only paste does not carry over formulas, so restore the
# "EARNED " helper formula in column G.
$ws.Cells.Item(37, 7).Formula = '=IF(ISBLANK([@EARNED]),"",[@EARNED])'

# Turn A37 into the bold "2024" year-separator label, matching the
# formatting already used by the "2022" label in A10.
$ws.Range("A10").Copy()
$ws.Range("A37").PasteSpecial(-4122)
$ws.Cells.Item(37, 1).Value = "'2024"

# ---------------------------------------------------------------------------
# Row 36 (12/31/2023): record the VL(1-0-0) holiday credit of 1 day, noting
# the 12/27/2023 reference date in the REMARKS column.
# ---------------------------------------------------------------------------
$ws.Cells.Item(36, 2).Value = "VL(1-0-0)"
$ws.Cells.Item(36, 4).Value = 1

$ws.Range("K31").Copy()
$ws.Range("K36").PasteSpecial(-4122)
$ws.Cells.Item(36, 11).Value = 45287

# ---------------------------------------------------------------------------
# Rows 33-35 (9/30/2023, 10/31/2023, 11/30/2023): record the 1.25 EARNED
# vacation-leave credit for each month.
# ---------------------------------------------------------------------------
$ws.Cells.Item(33, 3).Value = 1.25
$ws.Cells.Item(34, 3).Value = 1.25
$ws.Cells.Item(35, 3).Value = 1.25
